$d = $word.ActiveDocument

function Get-ParaByText($doc, $wanted) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $wanted) {
            return $p
        }
    }
    return $null
}

function Insert-ParaFragment($para, $innerXml) {
    $rng = $para.Range.Duplicate
    $rng.MoveEnd(1, -1)
    $xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

function Add-SpellSuffix($doc, $oldText, $suffix) {
    $p = Get-ParaByText $doc $oldText
    $inner = '<w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>' + $suffix + '</w:t></w:r><w:proofErr w:type="spellEnd"/>'
    Insert-ParaFragment $p $inner
}

# 1. Remove the _GoBack bookmark (currently after "WAT ER IN DE DATABASE MOET:")
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# 2. "Acteur" -> "Acteur actor.list"
Add-SpellSuffix $d "Acteur" "actor.list"

# 3. "Actrice" -> "Actrice actor.list"
Add-SpellSuffix $d "Actrice" "actor.list"

# 4. "Films" -> "Films actors.list"
Add-SpellSuffix $d "Films" "actors.list"

# 5. "Jaar film is uitgebracht" -> "Jaar film is uitgebracht actors.list"
Add-SpellSuffix $d "Jaar film is uitgebracht" "actors.list"

# 6. Re-add the _GoBack bookmark, now right after "Aantal stemmen film"
$p11 = Get-ParaByText $d "Aantal stemmen film"
Insert-ParaFragment $p11 '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
